$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name
$ws.Name = "Through 2022-05-25"

# Apply cell-level changes derived from the diff
$ws.Range("B1").Value = "May 2022 (through May 25)"
$ws.Range("B2").Value = 9
$ws.Range("AF6").Value = 2
$ws.Range("Q8").Value = 3
$ws.Range("L15").Value = 2
$ws.Range("A16").Value = "South Chicago"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 2
$ws.Range("E16").ClearContents()
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2
$ws.Range("O16").ClearContents()
$ws.Range("S16").ClearContents()
$ws.Range("U16").ClearContents()
$ws.Range("X16").Value = 1
$ws.Range("Y16").Value = 1
$ws.Range("AA16").Value = 1
$ws.Range("AC16").Value = 1
$ws.Range("AD16").Value = 1
$ws.Range("AE16").ClearContents()
$ws.Range("AF16").ClearContents()
$ws.Range("AJ16").ClearContents()
$ws.Range("AM16").ClearContents()
$ws.Range("AN16").Value = 2
$ws.Range("AO16").ClearContents()
$ws.Range("A17").Value = "Pullman"
$ws.Range("B17").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()
$ws.Range("H17").Value = 2
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("X17").ClearContents()
$ws.Range("Y17").ClearContents()
$ws.Range("Z17").ClearContents()
$ws.Range("AA17").ClearContents()
$ws.Range("AB17").ClearContents()
$ws.Range("AC17").ClearContents()
$ws.Range("AD17").ClearContents()
$ws.Range("AN17").Value = 1
$ws.Range("A18").Value = "Montclare"
$ws.Range("AN18").ClearContents()
$ws.Range("A19").Value = "Little Italy, UIC"
$ws.Range("B19").Value = 2
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 6
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 1
$ws.Range("S19").Value = 2
$ws.Range("U19").Value = 2
$ws.Range("Y19").Value = 2
$ws.Range("Z19").Value = 1
$ws.Range("AB19").Value = 2
$ws.Range("AE19").Value = 1
$ws.Range("AF19").Value = 1
$ws.Range("AJ19").Value = 1
$ws.Range("AM19").Value = 3
$ws.Range("AN19").Value = 1
$ws.Range("AO19").Value = 2
$ws.Range("L20").Value = 2
$ws.Range("G21").Value = 2
$ws.Range("AA23").Value = 4
$ws.Range("L24").Value = 1
$ws.Range("AF24").Value = 1
$ws.Range("B26").Value = 1
$ws.Range("B30").Value = 3
$ws.Range("V32").Value = 1
$ws.Range("L41").Value = 1
$ws.Range("AK45").Value = 1
$ws.Range("B48").Value = 1
